$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date" column (column B) demo data.
# Rows 3-5 previously shared the date "2024-08-25" and are updated together
# to "2024-10-25"; row 2's date moves from "2024-08-01" to "2024-09-30".
$ws.Range("B3").Value = "2024-10-25"
$ws.Range("B4").Value = "2024-10-25"
$ws.Range("B5").Value = "2024-10-25"
$ws.Range("B2").Value = "2024-09-30"

# Leave the selection on B3, matching the saved cursor position.
$ws.Range("B3").Select()
